# nowcasts_2025Q4_Nr1_Np1_Nj1.xlsx — add results from latest run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise the B:K nowcast/revision figures for the existing vintages (rows 2-11) ---
# Row 2
$ws.Range("B2").Value = [double]"0.29180857720257425"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Row 3
$ws.Range("B3").Value = [double]"0.29179467479371896"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = [double]"5.2377154102540875e-06"
$ws.Range("E3").Value = [double]"-4.32201504328257e-07"
$ws.Range("F3").Value = [double]"-5.1669404282718171e-08"
$ws.Range("G3").Value = [double]"-3.4655490781911241e-07"
$ws.Range("H3").Value = [double]"-2.9437304133139655e-08"
$ws.Range("I3").Value = [double]"-1.2545094737201008e-06"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = [double]"-1.6995802487596023e-05"

# Row 4
$ws.Range("B4").Value = [double]"0.29155745328959654"
$ws.Range("C4").Value = [double]"-2.7463756931422666e-05"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = [double]"1.0452752701318788e-06"
$ws.Range("F4").Value = [double]"-3.4059729041119177e-06"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = [double]"-3.8566844274809354e-06"
$ws.Range("I4").Value = [double]"-4.1974836158737086e-05"
$ws.Range("J4").Value = [double]"-3.9575207538764989e-05"
$ws.Range("K4").Value = [double]"-1.5737071573818096e-05"

# Row 5
$ws.Range("B5").Value = [double]"0.29249458244124582"
$ws.Range("C5").Value = [double]"2.4112953621296108e-05"
$ws.Range("D5").Value = [double]"0.00023288776501718126"
$ws.Range("E5").Value = [double]"5.5753262299166558e-05"
$ws.Range("F5").Value = [double]"0.00013116278447671199"
$ws.Range("G5").Value = [double]"-9.4624779553512817e-06"
$ws.Range("H5").Value = [double]"8.879533571393353e-06"
$ws.Range("I5").Value = [double]"-2.476905921594688e-06"
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = [double]"5.7288534141608505e-05"

# Row 6
$ws.Range("B6").Value = [double]"0.29332235918831101"
$ws.Range("C6").Value = [double]"0.0012667107939491129"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = [double]"-3.4231124439295003e-06"
$ws.Range("F6").Value = [double]"-4.5829996064742957e-06"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = [double]"5.584846949284453e-07"
$ws.Range("I6").Value = [double]"-2.1921735887829625e-05"
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = [double]"-0.00070702462358079865"

# Row 7
$ws.Range("B7").Value = [double]"0.29274577443827393"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = [double]"-0.00018646233140194075"
$ws.Range("E7").Value = [double]"1.4474639632467839e-07"
$ws.Range("F7").Value = [double]"-4.2867082765333545e-06"
$ws.Range("G7").Value = [double]"1.6043687739984064e-05"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = [double]"6.1074976007294311e-06"
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = [double]"-0.00037592543599962447"

# Row 8
$ws.Range("B8").Value = [double]"0.28841244425076279"
$ws.Range("C8").Value = [double]"-0.0039973312730613589"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = [double]"3.1572913856478985e-07"
$ws.Range("F8").Value = [double]"-1.4224908330645857e-05"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = [double]"-8.5087805718731297e-07"
$ws.Range("I8").Value = [double]"0.0001616428566976706"
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = [double]"-2.2617538632530021e-06"

# Row 9
$ws.Range("B9").Value = [double]"0.28717283336378457"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = [double]"-0.00055813080745866297"
$ws.Range("E9").Value = [double]"-3.9719309519645994e-05"
$ws.Range("F9").Value = [double]"-0.00020826765663065809"
$ws.Range("G9").Value = [double]"-3.4321360790166818e-05"
$ws.Range("H9").Value = [double]"-7.3853147883835908e-06"
$ws.Range("I9").Value = [double]"8.568795645293963e-06"
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = [double]"-1.0436329137830658e-05"

# Row 10
$ws.Range("B10").Value = [double]"0.30978013341989891"
$ws.Range("C10").Value = [double]"0.023201758456323756"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = [double]"1.2200087754986122e-05"
$ws.Range("F10").Value = [double]"1.9856845403149657e-05"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = [double]"1.3474028664034264e-05"
$ws.Range("I10").Value = [double]"-6.9880241850632641e-07"
$ws.Range("J10").Value = [double]"-0.002202230105198081"
$ws.Range("K10").Value = [double]"0.0008636547946184181"

# Row 11
$ws.Range("B11").Value = [double]"0.31762194009609546"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = [double]"0.0065872565240356973"
$ws.Range("E11").Value = [double]"-1.2805020467949557e-05"
$ws.Range("F11").Value = [double]"-7.0557259336400399e-05"
$ws.Range("G11").Value = [double]"0.00022491995958260811"
$ws.Range("H11").Value = [double]"-1.0963649442072172e-05"
$ws.Range("I11").Value = [double]"-0.00091800808960951436"
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = [double]"0.0022173454029685291"

# --- Append the new vintage row for 2025-08-30 ---
# (leading apostrophe forces text instead of Excel's date auto-detection;
#  re-apply the plain style used by the other rows so no quotePrefix sticks)
$ws.Range("A12").Formula = "'2025-08-30"
$ws.Range("A12").Style = $ws.Range("A11").Style
$ws.Range("B12").Value = [double]"0.27067844482950693"
$ws.Range("C12").Value = [double]"-0.034747241197966842"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = [double]"-1.7005954854300203e-06"
$ws.Range("F12").Value = [double]"-3.9978308232488955e-07"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = [double]"1.3142646085194413e-06"
$ws.Range("I12").Value = [double]"-0.00017226292559806961"
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = [double]"-0.003749760710910377"

# --- Widen column I slightly to match the refreshed layout ---
$ws.Columns.Item(9).ColumnWidth = [double]"15.3"
